$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Sheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 511.97778
$ws.Range("J17").Value = 298.35135
$ws.Range("L17").Value = 895.0540500000001
$ws.Range("N17").Value = -1231.05405
# Row 51
$ws.Range("H51").Value = 5362.727
$ws.Range("J51").Value = 7518
$ws.Range("L51").Value = 7518
$ws.Range("N51").Value = -8486
# Row 63
$ws.Range("H63").Value = 25333
$ws.Range("J63").Value = 25333
$ws.Range("L63").Value = 25333
$ws.Range("N63").Value = -26581
# Row 66
$ws.Range("H66").Value = 25333
$ws.Range("J66").Value = 25333
$ws.Range("L66").Value = 75999
$ws.Range("N66").Value = -82239
# Row 100
$ws.Range("H100").Value = 28573626
$ws.Range("I100").Value = 28573626
$ws.Range("K100").Value = 28573626
$ws.Range("M100").Value = -28573085
# Row 109
$ws.Range("H109").Value = 34657.895
$ws.Range("J109").Value = 34657.895
$ws.Range("L109").Value = 34657.895
$ws.Range("N109").Value = -37431.895
# Row 111
$ws.Range("H111").Value = 850.4666999999999
$ws.Range("I111").Value = 524.36365
$ws.Range("K111").Value = 1573.09095
$ws.Range("M111").Value = 1493.90905
# Row 112
$ws.Range("H112").Value = 1314.4642
$ws.Range("J112").Value = 1314.4642
$ws.Range("L112").Value = 3943.3926
$ws.Range("N112").Value = -6159.392599999999
# Row 121
$ws.Range("H121").Value = 2490
$ws.Range("J121").Value = 2490
$ws.Range("L121").Value = 7470
$ws.Range("N121").Value = -10964
# Row 129
$ws.Range("H129").Value = 829.8687
$ws.Range("I129").Value = 314.66666
$ws.Range("J129").Value = 863.10754
$ws.Range("K129").Value = 943.9999799999999
$ws.Range("L129").Value = 2589.32262
$ws.Range("M129").Value = 4056.00002
$ws.Range("N129").Value = -12589.32262
# Row 132
$ws.Range("H132").Value = 45461760
$ws.Range("I132").Value = 52638776
$ws.Range("J132").Value = 7333.3335
$ws.Range("K132").Value = 157916328
$ws.Range("L132").Value = 22000.0005
$ws.Range("M132").Value = -157913798
$ws.Range("N132").Value = -27060.0005

# ---- Sheet ARM ----
$ws = $wb.Sheets.Item("ARM")
# Row 64
$ws.Range("H64").Value = 26522.5
$ws.Range("J64").Value = 26522.5
$ws.Range("L64").Value = 26522.5
$ws.Range("N64").Value = -27018.5
# Row 67
$ws.Range("H67").Value = 26522.5
$ws.Range("J67").Value = 26522.5
$ws.Range("L67").Value = 26522.5
$ws.Range("N67").Value = -28238.5
# Row 102
$ws.Range("H102").Value = 2000
$ws.Range("I102").Value = 2000
$ws.Range("K102").Value = 2000
$ws.Range("M102").Value = -378

# ---- Sheet BSM ----
$ws = $wb.Sheets.Item("BSM")
# Row 40
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
# Row 62
$ws.Range("H62").Value = 36900
$ws.Range("J62").Value = 36900
$ws.Range("L62").Value = 36900
$ws.Range("N62").Value = -38272
# Row 65
$ws.Range("H65").Value = 36900
$ws.Range("J65").Value = 36900
$ws.Range("L65").Value = 110700
$ws.Range("N65").Value = -117564
# Row 94
$ws.Range("H94").Value = 641.9
$ws.Range("I94").Value = 633.625
$ws.Range("K94").Value = 633.625
$ws.Range("M94").Value = -182.625
# Row 96
$ws.Range("H96").Value = 4428
$ws.Range("I96").Value = 4428
$ws.Range("K96").Value = 4428
$ws.Range("M96").Value = -1682
# Row 134
$ws.Range("H134").Value = 3338.16
$ws.Range("I134").Value = 1707.7646
$ws.Range("J134").Value = 6802.75
$ws.Range("K134").Value = 5123.293799999999
$ws.Range("L134").Value = 20408.25
$ws.Range("M134").Value = -2588.293799999999
$ws.Range("N134").Value = -25478.25

# ---- Sheet CRP ----
$ws = $wb.Sheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 7018.5483
$ws.Range("I31").Value = 2881
$ws.Range("K31").Value = 2881
$ws.Range("M31").Value = -2586
# Row 34
$ws.Range("H34").Value = 7018.5483
$ws.Range("I34").Value = 2881
$ws.Range("K34").Value = 2881
$ws.Range("M34").Value = -2679
# Row 115
$ws.Range("H115").Value = 34799
$ws.Range("J115").Value = 34799
$ws.Range("L115").Value = 34799
$ws.Range("N115").Value = -37149
# Row 122
$ws.Range("H122").Value = 3678.4285
$ws.Range("J122").Value = 5000
$ws.Range("L122").Value = 15000
$ws.Range("N122").Value = -19900
# Row 139
$ws.Range("H139").Value = 38824.75
$ws.Range("J139").Value = 38824.75
$ws.Range("L139").Value = 38824.75
$ws.Range("N139").Value = -49104.75

# ---- Sheet CUL ----
$ws = $wb.Sheets.Item("CUL")
# Row 40
$ws.Range("H40").Value = 120.55556
$ws.Range("I40").Value = 106.42857
$ws.Range("J40").Value = 170
$ws.Range("K40").Value = 425.71428
$ws.Range("L40").Value = 680
$ws.Range("M40").Value = -356.71428
$ws.Range("N40").Value = -818
# Row 107
$ws.Range("H107").Value = 25050720
$ws.Range("I107").Value = 548.8333
$ws.Range("J107").Value = 35786508
$ws.Range("K107").Value = 1646.4999
$ws.Range("L107").Value = 107359524
$ws.Range("M107").Value = 273.5001
$ws.Range("N107").Value = -107363364
# Row 131
$ws.Range("H131").Value = 769.8
$ws.Range("I131").Value = 300
$ws.Range("J131").Value = 822
$ws.Range("K131").Value = 900
$ws.Range("L131").Value = 2466
$ws.Range("M131").Value = 4140
$ws.Range("N131").Value = -12546

# ---- Sheet GSM ----
$ws = $wb.Sheets.Item("GSM")
# Row 107
$ws.Range("H107").Value = 5848621
$ws.Range("J107").Value = 13889751
$ws.Range("L107").Value = 13889751
$ws.Range("N107").Value = -13893591

# ---- Sheet LTW ----
$ws = $wb.Sheets.Item("LTW")
# Row 81
$ws.Range("H81").Value = 53528.285
$ws.Range("J81").Value = 53528.285
$ws.Range("L81").Value = 53528.285
$ws.Range("N81").Value = -55524.285
# Row 84
$ws.Range("H84").Value = 53528.285
$ws.Range("J84").Value = 53528.285
$ws.Range("L84").Value = 160584.855
$ws.Range("N84").Value = -170568.855
# Row 92
$ws.Range("H92").Value = 29450
$ws.Range("J92").Value = 29450
$ws.Range("L92").Value = 29450
$ws.Range("N92").Value = -34442
# Row 101
$ws.Range("H101").Value = 44362
$ws.Range("J101").Value = 44362
$ws.Range("L101").Value = 44362
$ws.Range("N101").Value = -50852
# Row 132
$ws.Range("H132").Value = 5609.5483
$ws.Range("I132").Value = 3092.7144
$ws.Range("J132").Value = 7682.2354
$ws.Range("K132").Value = 9278.143199999999
$ws.Range("L132").Value = 23046.7062
$ws.Range("M132").Value = -6748.143199999999
$ws.Range("N132").Value = -28106.7062

# ---- Sheet WVR ----
$ws = $wb.Sheets.Item("WVR")
# Row 64
$ws.Range("H64").Value = 28900
$ws.Range("J64").Value = 28900
$ws.Range("L64").Value = 28900
$ws.Range("N64").Value = -29396
# Row 67
$ws.Range("H67").Value = 28900
$ws.Range("J67").Value = 28900
$ws.Range("L67").Value = 28900
$ws.Range("N67").Value = -30616
# Row 80
$ws.Range("H80").Value = 39790
$ws.Range("J80").Value = 39790
$ws.Range("L80").Value = 39790
$ws.Range("N80").Value = -41786
# Row 83
$ws.Range("H83").Value = 39790
$ws.Range("J83").Value = 39790
$ws.Range("L83").Value = 119370
$ws.Range("N83").Value = -129354
# Row 103
$ws.Range("H103").Value = 33534
$ws.Range("J103").Value = 33534
$ws.Range("L103").Value = 33534
$ws.Range("N103").Value = -35878
# Row 122
$ws.Range("H122").Value = 7340.1113
$ws.Range("I122").Value = 4942.4546
$ws.Range("J122").Value = 11107.857
$ws.Range("K122").Value = 14827.3638
$ws.Range("L122").Value = 33323.571
$ws.Range("M122").Value = -12377.3638
$ws.Range("N122").Value = -38223.571
# Row 128
$ws.Range("H128").Value = 41378.332
$ws.Range("J128").Value = 41378.332
$ws.Range("L128").Value = 41378.332
$ws.Range("N128").Value = -51338.332
